$d = $word.ActiveDocument

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# "git add <newfiletoaddname>" - split the way Word's spell checker would,
# flagging the inner word as a (misspelled) run of its own.
$para1 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>git add &lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>newfiletoaddname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&gt;</w:t></w:r></w:p>'

# git commit -m "commit message"  (curly/smart quotes, as used elsewhere in the doc)
$para2 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>git commit -m &#8220;commit message&#8221;</w:t></w:r></w:p>'

# git push origin master  (trailing space preserved)
$para3 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">git push origin master </w:t></w:r></w:p>'

function Insert-ParagraphXml($xmlFragment) {
    $doc = $word.ActiveDocument
    $lastPara = $doc.Paragraphs.Last
    $endPos = $lastPara.Range.End
    # Target the point just before the trailing paragraph mark so the new
    # paragraph is spliced in right after the current last paragraph,
    # instead of overwriting it / leaving a stray empty paragraph behind.
    $insertionPoint = $doc.Range($endPos - 1, $endPos - 1)
    [void]$insertionPoint.InsertXML($xmlFragment)
}

Insert-ParagraphXml ($pkgHeader + $para1 + $pkgFooter)
Insert-ParagraphXml ($pkgHeader + $para2 + $pkgFooter)
Insert-ParagraphXml ($pkgHeader + $para3 + $pkgFooter)
